$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header (row 1)
$ws.Range('A1').Value = 'Datos actualizados a 19 de Abril de 2020 a las 18:52'

# Row 4: Estados Unidos
$ws.Range('A4').Value = 'Estados Unidos'
$ws.Range('B4').Value = 742732
$ws.Range('C4').Value = 3940
$ws.Range('D4').Value = 68658
$ws.Range('E4').Value = 634414
$ws.Range('F4').Value = 13551
$ws.Range('G4').Value = 646
$ws.Range('H4').Value = 39660

# Row 7: Francia
$ws.Range('A7').Value = 'Francia'
$ws.Range('B7').Value = 152578
$ws.Range('C7').Value = 785
$ws.Range('D7').Value = 36570
$ws.Range('E7').Value = 96290
$ws.Range('F7').Value = 5744
$ws.Range('G7').Value = 395
$ws.Range('H7').Value = 19718

# Row 8: Alemania
$ws.Range('A8').Value = 'Alemania'
$ws.Range('B8').Value = 144387
$ws.Range('C8').Value = 663
$ws.Range('D8').Value = 88000
$ws.Range('E8').Value = 51840
$ws.Range('F8').Value = 2922
$ws.Range('G8').Value = 9
$ws.Range('H8').Value = 4547

# Row 20: India
$ws.Range('A20').Value = 'India'
$ws.Range('B20').Value = 17615
$ws.Range('C20').Value = 1250
$ws.Range('D20').Value = 2769
$ws.Range('E20').Value = 14290
$ws.Range('F20').Value = 0
$ws.Range('G20').Value = 35
$ws.Range('H20').Value = 556

# Row 26: Japon
$ws.Range('A26').Value = 'Japon'
$ws.Range('B26').Value = 10797
$ws.Range('C26').Value = 501
$ws.Range('D26').Value = 1159
$ws.Range('E26').Value = 9402
$ws.Range('F26').Value = 217
$ws.Range('G26').Value = 14
$ws.Range('H26').Value = 236

# Row 38: Chequia
$ws.Range('A38').Value = 'Chequia'
$ws.Range('B38').Value = 6701
$ws.Range('C38').Value = 95
$ws.Range('D38').Value = 1298
$ws.Range('E38').Value = 5217
$ws.Range('F38').Value = 84
$ws.Range('G38').Value = 5
$ws.Range('H38').Value = 186

# Row 52: Luxemburgo
$ws.Range('A52').Value = 'Luxemburgo'
$ws.Range('B52').Value = 3550
$ws.Range('C52').Value = 13
$ws.Range('D52').Value = 627
$ws.Range('E52').Value = 2850
$ws.Range('F52').Value = 31
$ws.Range('G52').Value = 1
$ws.Range('H52').Value = 73

# Row 80: Cuba
$ws.Range('A80').Value = 'Cuba'
$ws.Range('B80').Value = 1035
$ws.Range('C80').Value = 49
$ws.Range('D80').Value = 255
$ws.Range('E80').Value = 746
$ws.Range('F80').Value = 16
$ws.Range('G80').Value = 2
$ws.Range('H80').Value = 34

# Row 81: Hong Kong
$ws.Range('A81').Value = 'Hong Kong'
$ws.Range('B81').Value = 1026
$ws.Range('C81').Value = 2
$ws.Range('D81').Value = 602
$ws.Range('E81').Value = 420
$ws.Range('F81').Value = 8
$ws.Range('G81').Value = 0
$ws.Range('H81').Value = 4

# Row 82: Camerun
$ws.Range('A82').Value = 'Camerun'
$ws.Range('B82').Value = 1017
$ws.Range('C82').Value = 0
$ws.Range('D82').Value = 305
$ws.Range('E82').Value = 670
$ws.Range('F82').Value = 33
$ws.Range('G82').Value = 0
$ws.Range('H82').Value = 42

# Row 83: Afganistan
$ws.Range('A83').Value = 'Afganistan'
$ws.Range('B83').Value = 996
$ws.Range('C83').Value = 63
$ws.Range('D83').Value = 131
$ws.Range('E83').Value = 832
$ws.Range('F83').Value = 7
$ws.Range('G83').Value = 3
$ws.Range('H83').Value = 33

# Row 96: Burkina Faso
$ws.Range('A96').Value = 'Burkina Faso'
$ws.Range('B96').Value = 576
$ws.Range('C96').Value = 11
$ws.Range('D96').Value = 338
$ws.Range('E96').Value = 202
$ws.Range('F96').Value = 0
$ws.Range('G96').Value = 0
$ws.Range('H96').Value = 36

# Row 109: Reunion
$ws.Range('A109').Value = 'Reunion'
$ws.Range('B109').Value = 408
$ws.Range('C109').Value = 1
$ws.Range('D109').Value = 237
$ws.Range('E109').Value = 171
$ws.Range('F109').Value = 4
$ws.Range('G109').Value = 0
$ws.Range('H109').Value = 0

# Row 111: Senegal
$ws.Range('A111').Value = 'Senegal'
$ws.Range('B111').Value = 367
$ws.Range('C111').Value = 17
$ws.Range('D111').Value = 220
$ws.Range('E111').Value = 144
$ws.Range('F111').Value = 1
$ws.Range('G111').Value = 0
$ws.Range('H111').Value = 3

# Row 116: Sri Lanka
$ws.Range('A116').Value = 'Sri Lanka'
$ws.Range('B116').Value = 271
$ws.Range('C116').Value = 17
$ws.Range('D116').Value = 96
$ws.Range('E116').Value = 168
$ws.Range('F116').Value = 1
$ws.Range('G116').Value = 0
$ws.Range('H116').Value = 7

# Row 117: Mayotte
$ws.Range('A117').Value = 'Mayotte'
$ws.Range('B117').Value = 271
$ws.Range('C117').Value = 17
$ws.Range('D117').Value = 117
$ws.Range('E117').Value = 150
$ws.Range('F117').Value = 5
$ws.Range('G117').Value = 0
$ws.Range('H117').Value = 4

# Row 118: Kenia
$ws.Range('A118').Value = 'Kenia'
$ws.Range('B118').Value = 270
$ws.Range('C118').Value = 8
$ws.Range('D118').Value = 67
$ws.Range('E118').Value = 189
$ws.Range('F118').Value = 2
$ws.Range('G118').Value = 2
$ws.Range('H118').Value = 14

# Row 119: Vietnam
$ws.Range('A119').Value = 'Vietnam'
$ws.Range('B119').Value = 268
$ws.Range('C119').Value = 0
$ws.Range('D119').Value = 203
$ws.Range('E119').Value = 65
$ws.Range('F119').Value = 8
$ws.Range('G119').Value = 0
$ws.Range('H119').Value = 0

# Row 120: Guatemala
$ws.Range('A120').Value = 'Guatemala'
$ws.Range('B120').Value = 257
$ws.Range('C120').Value = 22
$ws.Range('D120').Value = 21
$ws.Range('E120').Value = 229
$ws.Range('F120').Value = 3
$ws.Range('G120').Value = 0
$ws.Range('H120').Value = 7

# Row 122: Mali
$ws.Range('A122').Value = 'Mali'
$ws.Range('B122').Value = 224
$ws.Range('C122').Value = 8
$ws.Range('D122').Value = 42
$ws.Range('E122').Value = 168
$ws.Range('F122').Value = 0
$ws.Range('G122').Value = 1
$ws.Range('H122').Value = 14

# Row 152: Guyana
$ws.Range('A152').Value = 'Guyana'
$ws.Range('B152').Value = 63
$ws.Range('C152').Value = 0
$ws.Range('D152').Value = 9
$ws.Range('E152').Value = 47
$ws.Range('F152').Value = 4
$ws.Range('G152').Value = 1
$ws.Range('H152').Value = 7
